$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountCreation")

# Row 2 holds the details for a generated test account. Update the
# first/last name, email, password and phone number for the new account,
# preserving each cell's existing number formatting (e.g. the phone
# number cell is stored with a quote-prefix / text format).
$ws.Range("J2").Copy($ws.Range("Z99"))

$ws.Range("A2").Value = "Marina"
$ws.Range("B2").Value = "Avery"
$ws.Range("C2").Value = "testaccount25@email.com"
$ws.Range("D2").Value = "testaccount25"
$ws.Range("J2").Value = "8173678441"

$ws.Range("Z99").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("Z99").Clear()
$excel.CutCopyMode = $false
